# Updates cryptos list (prices / volume / two row-swaps) per the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Helper: writes $val into $ref as literal TEXT (not a number/date),
# then clears the temporary "@" number-format so the cell keeps the
# workbook default style (matches the source file, which never sets
# a style on these data cells).
function Set-TextValue($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.117.35"
$ws.Range("E2").Value = "  +0.02%  "
# Row 3
$ws.Range("D3").Value = "1.801.67"
$ws.Range("E3").Value = "  +2.32%  "
# Row 4
Set-TextValue $ws "D4" "1.007"
$ws.Range("E4").Value = "  +0.53%  "
# Row 5
Set-TextValue $ws "D5" "338.04"
$ws.Range("E5").Value = "  +0.77%  "
# Row 6
$ws.Range("E6").Value = "  +0.44%  "
# Row 7
Set-TextValue $ws "D7" "0.4736"
$ws.Range("E7").Value = "  +25.14%  "
# Row 8
Set-TextValue $ws "D8" "0.3824"
$ws.Range("E8").Value = "  +13.71%  "
# Row 9
$ws.Range("E9").Value = "  -0.47%  "
# Row 10
Set-TextValue $ws "D10" "0.07706"
$ws.Range("E10").Value = "  +6.74%  "
# Row 11
Set-TextValue $ws "D11" "1.143"
$ws.Range("E11").Value = "  +1.23%  "
# Row 12
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D12" "1.006"
$ws.Range("E12").Value = "  +0.63%  "
# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws "D13" "22.52"
$ws.Range("E13").Value = "  -0.65%  "
# Row 14
Set-TextValue $ws "D14" "6.427"
$ws.Range("E14").Value = "  +3.61%  "
# Row 15
Set-TextValue $ws "D15" "7.366"
$ws.Range("E15").Value = "  +2.04%  "
# Row 16
$ws.Range("D16").Value = "1.799.72"
$ws.Range("E16").Value = "  +2.36%  "
# Row 17
Set-TextValue $ws "D17" "0.00001097"
$ws.Range("E17").Value = "  +3.85%  "
# Row 18
Set-TextValue $ws "D18" "0.06745"
$ws.Range("E18").Value = "  +2.35%  "
# Row 19
Set-TextValue $ws "D19" "82.13"
$ws.Range("E19").Value = "  +1.69%  "
# Row 20
Set-TextValue $ws "D20" "1.000"
$ws.Range("E20").Value = "  +0.23%  "
# Row 21
Set-TextValue $ws "D21" "17.46"
$ws.Range("E21").Value = "  +2.51%  "
# Row 22
Set-TextValue $ws "D22" "6.448"
$ws.Range("E22").Value = "  +2.70%  "
# Row 23
$ws.Range("D23").Value = "28.154.54"
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
Set-TextValue $ws "D24" "11.97"
$ws.Range("E24").Value = "  +2.19%  "
# Row 25
Set-TextValue $ws "D25" "2.394"
$ws.Range("E25").Value = "  +0.35%  "
# Row 26
Set-TextValue $ws "D26" "20.93"
$ws.Range("E26").Value = "  +5.19%  "
# Row 27
Set-TextValue $ws "D27" "2.408"
$ws.Range("E27").Value = "  +2.51%  "
# Row 28
Set-TextValue $ws "D28" "151.51"
$ws.Range("E28").Value = "  -1.33%  "
# Row 29
$ws.Range("D29").Value = "2.010.23"
$ws.Range("E29").Value = "  +2.56%  "
# Row 30
Set-TextValue $ws "D30" "134.02"
$ws.Range("E30").Value = "  +1.54%  "
# Row 31
Set-TextValue $ws "D31" "1.272"
$ws.Range("E31").Value = "  +0.01%  "
# Row 32
Set-TextValue $ws "D32" "4.067"
$ws.Range("E32").Value = "  +1.29%  "
# Row 33
Set-TextValue $ws "D33" "0.09769"
$ws.Range("E33").Value = "  +10.88%  "
# Row 34
Set-TextValue $ws "D34" "5.948"
$ws.Range("E34").Value = "  +1.91%  "
# Row 35
$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D35" "0.2244"
$ws.Range("E35").Value = "  +5.86%  "
# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D36" "0.02383"
$ws.Range("E36").Value = "  +1.53%  "
# Row 37
Set-TextValue $ws "D37" "12.19"
$ws.Range("E37").Value = "  -1.06%  "
# Row 38
Set-TextValue $ws "D38" "0.06387"
$ws.Range("E38").Value = "  +2.75%  "
# Row 39
Set-TextValue $ws "D39" "0.6698"
$ws.Range("E39").Value = "  +0.68%  "
# Row 40
$ws.Range("E40").Value = "  +1.14%  "
# Row 41
Set-TextValue $ws "D41" "1.493"
$ws.Range("E41").Value = "  +2.37%  "
# Row 42
Set-TextValue $ws "D42" "1.224"
$ws.Range("E42").Value = "  +0.68%  "
# Row 43
Set-TextValue $ws "D43" "8.096"
$ws.Range("E43").Value = "  +0.50%  "
# Row 44
Set-TextValue $ws "D44" "14.20"
$ws.Range("E44").Value = "  +3.22%  "
# Row 45
Set-TextValue $ws "D45" "1.001"
$ws.Range("E45").Value = "  +0.46%  "
# Row 46
Set-TextValue $ws "D46" "0.6173"
$ws.Range("E46").Value = "  +1.71%  "
# Row 47
Set-TextValue $ws "D47" "3.873"
$ws.Range("E47").Value = "  +1.58%  "
# Row 48
Set-TextValue $ws "D48" "129.67"
$ws.Range("E48").Value = "  -0.23%  "
# Row 49
Set-TextValue $ws "D49" "2.057"
$ws.Range("E49").Value = "  +1.81%  "
# Row 50
Set-TextValue $ws "D50" "1.176"
$ws.Range("E50").Value = "  -1.05%  "
# Row 51
Set-TextValue $ws "D51" "0.07125"
$ws.Range("E51").Value = "  -1.25%  "
